$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Tuesday, Jan 10"
$ws.Range("C23").Value = "12:05 AM"
$ws.Range("D23").Value = "PQ7552"
$ws.Range("E23").Value = "Sharm el-Sheikh"
$ws.Range("F23").Value = "(SSH)"
$ws.Range("G23").Value = "SmartWings "
$ws.Range("H23").Value = "B738"
$ws.Range("I23").Value = "(UR-SQQ)"
$ws.Range("J23").Value = "12:23 AM"
$ws.Range("L23").Value = "0 hours, 18 minutes"

# Row 24
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Tuesday, Jan 10"
$ws.Range("C24").Value = "5:40 AM"
$ws.Range("D24").Value = "LO3804"
$ws.Range("E24").Value = "Warsaw"
$ws.Range("F24").Value = "(WAW)"
$ws.Range("G24").Value = "LOT "
$ws.Range("H24").Value = "E190"
$ws.Range("I24").Value = "(SP-LMB)"
$ws.Range("J24").Value = "5:44 AM"
$ws.Range("L24").Value = "0 hours, 4 minutes"

# Row 25
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Tuesday, Jan 10"
$ws.Range("C25").Value = "10:00 AM"
$ws.Range("D25").Value = "PQ2025"
$ws.Range("E25").Value = "Antalya"
$ws.Range("F25").Value = "(AYT)"
$ws.Range("G25").Value = "SkyUp Airlines "
$ws.Range("H25").Value = "B738"
$ws.Range("I25").Value = "(UR-SQH)"
$ws.Range("J25").Value = "2:28 PM"
$ws.Range("L25").Value = "4 hours, 28 minutes"

# Row 26
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Tuesday, Jan 10"
$ws.Range("C26").Value = "10:05 AM"
$ws.Range("D26").Value = "FR2137"
$ws.Range("E26").Value = "London"
$ws.Range("F26").Value = "(STN)"
$ws.Range("G26").Value = "Ryanair "
$ws.Range("H26").Value = "B738"
$ws.Range("I26").Value = "(EI-DHE)"
$ws.Range("J26").Value = "10:36 AM"
$ws.Range("L26").Value = "0 hours, 31 minutes"

# Row 27
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "Tuesday, Jan 10"
$ws.Range("C27").Value = "11:35 AM"
$ws.Range("D27").Value = "UNKNOWN"
$ws.Range("E27").Value = "Farnborough"
$ws.Range("F27").Value = "(FAB)"
$ws.Range("G27").Value = "Air X Charter "
$ws.Range("H27").Value = "E35L"
$ws.Range("I27").Value = "(9H-JPC)"
$ws.Range("J27").Value = "11:58 AM"
$ws.Range("L27").Value = "0 hours, 23 minutes"

# Row 28
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Tuesday, Jan 10"
$ws.Range("C28").Value = "11:50 AM"
$ws.Range("D28").Value = "W61650"
$ws.Range("E28").Value = "Eindhoven"
$ws.Range("F28").Value = "(EIN)"
$ws.Range("G28").Value = "Wizz Air "
$ws.Range("H28").Value = "A320"
$ws.Range("I28").Value = "(HA-LYH)"
$ws.Range("J28").Value = "11:54 AM"
$ws.Range("L28").Value = "0 hours, 4 minutes"

# Row 29
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "Tuesday, Jan 10"
$ws.Range("C29").Value = "12:00 PM"
$ws.Range("D29").Value = "LO3810"
$ws.Range("E29").Value = "Warsaw"
$ws.Range("F29").Value = "(WAW)"
$ws.Range("G29").Value = "LOT "
$ws.Range("H29").Value = "E75S"
$ws.Range("I29").Value = "(SP-LIK)"
$ws.Range("J29").Value = "12:25 PM"
$ws.Range("L29").Value = "0 hours, 25 minutes"

# Row 30
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "Tuesday, Jan 10"
$ws.Range("C30").Value = "12:00 PM"
$ws.Range("D30").Value = "SK7182"
$ws.Range("E30").Value = "Oslo"
$ws.Range("F30").Value = "(OSL)"
$ws.Range("G30").Value = "SAS "
$ws.Range("H30").Value = "B737"
$ws.Range("I30").Value = "(LN-RPJ)"
$ws.Range("J30").Value = "12:17 PM"
$ws.Range("L30").Value = "0 hours, 17 minutes"

# Row 31
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "Tuesday, Jan 10"
$ws.Range("C31").Value = "2:55 PM"
$ws.Range("D31").Value = "LO3802"
$ws.Range("E31").Value = "Warsaw"
$ws.Range("F31").Value = "(WAW)"
$ws.Range("G31").Value = "LOT "
$ws.Range("H31").Value = "E190"
$ws.Range("I31").Value = "(SP-LMD)"
$ws.Range("J31").Value = "3:03 PM"
$ws.Range("L31").Value = "0 hours, 8 minutes"

# Row 32
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "Tuesday, Jan 10"
$ws.Range("C32").Value = "4:45 PM"
$ws.Range("D32").Value = "LO3808"
$ws.Range("E32").Value = "Warsaw"
$ws.Range("F32").Value = "(WAW)"
$ws.Range("G32").Value = "LOT "
$ws.Range("H32").Value = "E170"
$ws.Range("I32").Value = "(SP-LDH)"
$ws.Range("J32").Value = "4:40 PM"
$ws.Range("L32").Value = "0 hours, -5 minutes"

# Row 33
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "Tuesday, Jan 10"
$ws.Range("C33").Value = "6:20 PM"
$ws.Range("D33").Value = "FR3473"
$ws.Range("E33").Value = "London"
$ws.Range("F33").Value = "(LTN)"
$ws.Range("G33").Value = "Ryanair "
$ws.Range("H33").Value = "B738"
$ws.Range("I33").Value = "(EI-EMR)"
$ws.Range("J33").Value = "6:41 PM"
$ws.Range("L33").Value = "0 hours, 21 minutes"

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = "Tuesday, Jan 10"
$ws.Range("C34").Value = "8:15 PM"
$ws.Range("D34").Value = "RK3203"
$ws.Range("E34").Value = "Manchester"
$ws.Range("F34").Value = "(MAN)"
$ws.Range("G34").Value = "Ryanair "
$ws.Range("H34").Value = "B738"
$ws.Range("I34").Value = "(G-RUKH)"
$ws.Range("J34").Value = "8:48 PM"
$ws.Range("L34").Value = "0 hours, 33 minutes"

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = "Tuesday, Jan 10"
$ws.Range("C35").Value = "9:50 PM"
$ws.Range("D35").Value = "W95154"
$ws.Range("E35").Value = "London"
$ws.Range("F35").Value = "(LTN)"
$ws.Range("G35").Value = "Wizz Air "
$ws.Range("H35").Value = "A321"
$ws.Range("I35").Value = "(G-WUKG)"
$ws.Range("J35").Value = "10:02 PM"
$ws.Range("L35").Value = "0 hours, 12 minutes"

Write-Host "Added rows 23-35"